$d = $word.ActiveDocument

function Set-ParagraphText($paragraph, $newText) {
    $r = $paragraph.Range
    # Exclude the trailing paragraph mark so we only replace the run text,
    # keeping the paragraph's own pPr/style intact.
    $r.End = $r.End - 1

    # Force a real structural mutation (collapsing the many single-word
    # runs down to one run) even when the final text content is identical
    # to the concatenation of the existing runs: write a transient
    # placeholder first, then the real text. (A no-op replace - new text
    # equal to the existing concatenated text - is otherwise skipped.)
    $r.Text = "zzPLACEHOLDERzz"
    $r2 = $paragraph.Range
    $r2.End = $r2.End - 1
    $r2.Text = $newText
}

Set-ParagraphText $d.Paragraphs.Item(1) "Answers: Introduction to complex numbers"
Set-ParagraphText $d.Paragraphs.Item(2) "Tom Coleman"
Set-ParagraphText $d.Paragraphs.Item(4) "Answers to questions relating to the guide on introduction to complex numbers."
